$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 15 with job listing data (Informatica Cloud Data Governance Catalog Specialist)
$ws.Range("A15").Value = 'Informatica Cloud Data Governance Catalog Specialist'
$ws.Range("B15").Value = 'https://www.dice.com/job-detail/6f16b244-a92c-425a-9007-b06701ccabcf?filters.employmentType=CONTRACTS%7CTHIRD_PARTY&filters.postedDate=ONE&page=2&location=United+States&longitude=-106.5348379&latitude=38.7945952&locationPrecision=Country&countryCode=US&q=Golang'
$ws.Range("C15").Value = 'Torrance, California'
$ws.Range("D15").Value = 'Contract'
$ws.Range("E15").Value = 'USD 53.00 - 55.00 per hour'
$ws.Range("F15").Value = 'Robert Half'
$ws.Range("G15").Value = 'Description We are looking for an experienced Informatica Cloud Data Governance Catalog Specialist to join our team in Southern California. This position involves working on-site four days per week and offers a long-term contract opportunity. The ideal candidate will have a strong background in data governance, analytics, and business intelligence tools, coupled with a proactive approach to problem-solving and collaboration. Responsibilities: Create catalog quality reports to monitor and enhance data governance metrics across domains and sub-domains. Develop and showcase data governance dashboards tailored to different user roles, including Data Owners, Stewards, Engineers, and Privacy Officers. Collaborate with business and IT teams, including data stewards, catalog architects, and platform owners, to implement governance solutions. Execute profiling, sampling, and scanner setups using Informatica tools to ensure data quality. Apply expertise in metadata management, data modeling, and large-scale data analysis to support governance initiatives. Design and implement both traditional relational and modern big-data architectures based on organizational requirements. Utilize business intelligence tools such as Power BI and Tableau to create actionable insights and reports. Define compliance procedures and produce audit reports to meet regulatory requirements. Establish and support governance councils and operational frameworks using data catalog tools. Facilitate metadata ingestion and ensure adherence to data security and quality standards. Requirements Proficiency with tools like Informatica Cloud Data Governance Catalog and Cloud Data Quality. Hands-on experience in data modeling, metadata management, and large-scale data analysis. Familiarity with Collibra, Alation, and Glue Data Catalog. Strong understanding of entity-relationship modeling and data security practices. Expertise in business intelligence technologies such as Power BI and Tableau. Exceptional communication and presentation skills to effectively convey technical concepts. Analytical mindset with proven problem-solving abilities. Ability to work collaboratively as part of a team and build strong relationships with stakeholders. Technology Doesn''t Change the World, People Do. Robert Half is the world''s first and largest specialized talent solutions firm that connects highly qualified job seekers to opportunities at great companies. We offer contract, temporary and permanent placement solutions for finance and accounting, technology, marketing and creative, legal, and administrative and customer support roles. Robert Half works to put you in the best position to succeed. We provide access to top jobs, competitive compensation and benefits, and free online training. Stay on top of every opportunity - whenever you choose - even on the go. Download the Robert Half app and get 1-tap apply, notifications of AI-matched jobs, and much more. All applicants applying for U.S. job openings must be legally authorized to work in the United States. Benefits are available to contract/temporary professionals, including medical, vision, dental, and life and disability insurance. Hired contract/temporary professionals are also eligible to enroll in our company 401(k) plan. Visit roberthalf.gobenefits.net for more information. 2025 Robert Half. An Equal Opportunity Employer. M/F/Disability/Veterans. By clicking "Apply Now," you''re agreeing to Robert Half''s Terms of Use and Privacy Notice.'
$ws.Range("H15").Value = '24.4%%'
$ws.Range("I15").Value = 'Contract Independent | Contract W2 | On-site | USD $53.00 - 55.00 per hour'
